# Archives of past results
# Rename the scratch "Sheet5" to hold an archived CPUE~Landings regression
# summary table, replacing the old gradebook placeholder data with a table
# of regression results (independents / p-value / r2 / r2 adj / F stat /
# w-lag note) plus a list of the other model runs being archived.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet5")
$ws.Name = "TotalCPUE~Landings"

# Wipe the old leftover gradebook content before laying down the new table.
$ws.Cells.Clear()

# Row labels for the archived-run list (column A, rows 4-8) are entered
# first, followed by the table header row, the "N" flag, and finally the
# two B90 labels - matching how the sheet was actually built up.
$ws.Range("A4").Value = "E98"
$ws.Range("A5").Value = "E99"
$ws.Range("A6").Value = "P88"
$ws.Range("A7").Value = "T06"
$ws.Range("A8").Value = "T38"

$ws.Range("A1").Value = "Independents"
$ws.Range("B1").Value = "P-value"
$ws.Range("C1").Value = "r2"
$ws.Range("D1").Value = "r2 adj"
$ws.Range("E1").Value = "F stat"
$ws.Range("F1").Value = "w/ lag"

$ws.Range("F2").Value = "N"

$ws.Range("A2").Value = "B90 (1980)"
$ws.Range("A3").Value = "B90 (2003)"

$ws.Range("B2").Value = 0.056750000000000002
$ws.Range("C2").Value = 0.32329999999999998
$ws.Range("D2").Value = 0.2712
$ws.Range("E2").Value = 6.2110000000000003

# Column A is widened to fit the longest label ("Independents"/"B90 (1980)").
$ws.Columns.Item(1).ColumnWidth = 12.7

# Leave the selection where the author left it.
$null = $ws.Range("B3").Select()
